$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F8").Value = 1630
$ws1.Range("F12").Value = 2645
$ws1.Range("F14").Value = 1497
$ws1.Range("F15").Value = 6980
$ws1.Range("F19").Value = 5197
$ws1.Range("F22").Value = 219
$ws1.Range("F28").Value = 19
$ws1.Range("F29").Value = 174
$ws1.Range("F32").Value = 1125

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F8").Value = 207

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F9").Value = 1630
$ws4.Range("F14").Value = 2645
$ws4.Range("F15").Value = 1497
$ws4.Range("F16").Value = 207
$ws4.Range("F23").Value = 5197
$ws4.Range("F27").Value = 219
$ws4.Range("F35").Value = 19
$ws4.Range("F36").Value = 174
$ws4.Range("F39").Value = 1125
